# Apply updated cryptocurrency price/volume data to Sheet1 (A1:E51 table)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.703.33"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.633.69"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.59"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").Value = "1.859.72"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "1.606.92"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "26.686.43"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.57"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.71"
$ws.Range("E19").Value = "  +8.13%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.38"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.15"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.90"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  +4.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "1.225.18"
$ws.Range("E36").Value = "  +5.37%  "
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.35"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "1.768.16"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.31"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("E50").Value = "  +4.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("E51").Value = "  -0.18%  "
